# Update "想去人数" (F column) counts by +1 on both the "展览" and "全部类型"
# sheets, which carry duplicate copies of the same exhibition listing.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 112
    $ws.Range("F3").Value = 17
    $ws.Range("F4").Value = 969
}
